$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value2 = "₹5.14 Lakh"
$ws.Range("G1").Value2 = "EMIs from ₹10,044/month"
$ws.Range("B3").Value2 = "2022 Maruti Vitara Brezza ZXI PLUS"
$ws.Range("C3").Value2 = "20,059 KM"
$ws.Range("D3").Value2 = "MANUAL"
$ws.Range("E3").Value2 = "TOP MODEL"
$ws.Range("F3").Value2 = "₹8.70 Lakh"
$ws.Range("G3").Value2 = "EMIs from ₹16,560/month"
$ws.Range("B5").Value2 = "2019 Maruti Swift ZXI PLUS AMT"
$ws.Range("C5").Value2 = "36,023 KM"
$ws.Range("D5").Value2 = "AUTOMATIC"
$ws.Range("E5").Value2 = "TOP MODEL"
$ws.Range("F5").Value2 = "₹6.31 Lakh"
$ws.Range("G5").Value2 = "EMIs from ₹12,336/month"
$ws.Range("B6").Value2 = "2018 Maruti Celerio VXI"
$ws.Range("C6").Value2 = "52,238 KM"
$ws.Range("D6").Value2 = "MANUAL"
$ws.Range("E6").Value2 = "100% TYRE LIFE REMAINING"
$ws.Range("F6").Value2 = "₹3.98 Lakh"
$ws.Range("G6").Value2 = "EMIs from ₹7,781/month"
$ws.Range("B7").Value2 = "2012 Maruti Wagon R 1.0 VXI"
$ws.Range("C7").Value2 = "30,121 KM"
$ws.Range("D7").Value2 = "MANUAL"
$ws.Range("E7").Value2 = "TOP MODEL"
$ws.Range("F7").Value2 = "₹2.07 Lakh"
$ws.Range("G7").Value2 = "EMIs from ₹6,875/month"
$ws.Range("B8").Value2 = "2022 Maruti Celerio VXI CNG"
$ws.Range("C8").Value2 = "6,790 KM"
$ws.Range("D8").Value2 = "MANUAL"
$ws.Range("E8").Value2 = "REGULARLY SERVICED"
$ws.Range("F8").Value2 = "₹6.20 Lakh"
$ws.Range("G8").Value2 = "EMIs from ₹12,121/month"
$ws.Range("B9").Value2 = "2012 Maruti Swift Dzire VXI"
$ws.Range("C9").Value2 = "22,466 KM"
$ws.Range("D9").Value2 = "MANUAL"
$ws.Range("E9").Value2 = "REGULARLY SERVICED"
$ws.Range("F9").Value2 = "₹2.91 Lakh"
$ws.Range("G9").Value2 = "EMIs from ₹9,665/month"
$ws.Range("B10").Value2 = "2016 Maruti Baleno ZETA PETROL 1.2"
$ws.Range("C10").Value2 = "96,466 KM"
$ws.Range("D10").Value2 = "MANUAL"
$ws.Range("E10").Value2 = "ALLOY WHEELS"
$ws.Range("F10").Value2 = "₹4.59 Lakh"
$ws.Range("G10").Value2 = "EMIs from ₹8,974/month"
$ws.Range("B11").Value2 = "2022 Maruti Vitara Brezza VXI AT SHVS"
$ws.Range("C11").Value2 = "1,402 KM"
$ws.Range("D11").Value2 = "AUTOMATIC"
$ws.Range("E11").Value2 = "STANDARD SAFETY FEATURES"
$ws.Range("F11").Value2 = "₹10.64 Lakh"
$ws.Range("G11").Value2 = "EMIs from ₹20,252/month"
$ws.Range("B12").Value2 = "2018 Maruti Dzire VXI AMT"
$ws.Range("C12").Value2 = "45,657 KM"
$ws.Range("D12").Value2 = "AUTOMATIC"
$ws.Range("E12").Value2 = "REGULARLY SERVICED"
$ws.Range("F12").Value2 = "₹5.24 Lakh"
$ws.Range("G12").Value2 = "EMIs from ₹10,244/month"
$ws.Range("B14").Value2 = "2017 Maruti Swift ZXI"
$ws.Range("C14").Value2 = "63,367 KM"
$ws.Range("D14").Value2 = "MANUAL"
$ws.Range("E14").Value2 = "TOP MODEL"
$ws.Range("F14").Value2 = "₹5.43 Lakh"
$ws.Range("G14").Value2 = "EMIs from ₹10,616/month"
$ws.Range("B15").Value2 = "2022 Maruti Swift ZXI PLUS"
$ws.Range("C15").Value2 = "24,626 KM"
$ws.Range("D15").Value2 = "MANUAL"
$ws.Range("E15").Value2 = "TOP MODEL"
$ws.Range("F15").Value2 = "₹8.01 Lakh"
$ws.Range("G15").Value2 = "EMIs from ₹15,246/month"
$ws.Range("B16").Value2 = "2010 Maruti Wagon R 1.0 VXI"
$ws.Range("C16").Value2 = "84,625 KM"
$ws.Range("D16").Value2 = "MANUAL"
$ws.Range("E16").Value2 = "TOP MODEL"
$ws.Range("F16").Value2 = "₹1.42 Lakh"
$ws.Range("G16").Value2 = "EMIs from ₹12,617/month"
$ws.Range("B17").Value2 = "2022 Maruti Baleno ZETA PETROL 1.2"
$ws.Range("C17").Value2 = "15,515 KM"
$ws.Range("D17").Value2 = "MANUAL"
$ws.Range("E17").Value2 = "ALLOY WHEELS"
$ws.Range("F17").Value2 = "₹8.61 Lakh"
$ws.Range("G17").Value2 = "EMIs from ₹16,388/month"
$ws.Range("B18").Value2 = "2022 Maruti S PRESSO VXI (O) CNG"
$ws.Range("C18").Value2 = "21,240 KM"
$ws.Range("D18").Value2 = "MANUAL"
$ws.Range("E18").Value2 = "REGULARLY SERVICED"
$ws.Range("F18").Value2 = "₹4.75 Lakh"
$ws.Range("G18").Value2 = "EMIs from ₹9,286/month"
$ws.Range("B19").Value2 = "2019 Maruti Baleno ZETA PETROL 1.2"
$ws.Range("C19").Value2 = "30,890 KM"
$ws.Range("D19").Value2 = "MANUAL"
$ws.Range("E19").Value2 = "ALLOY WHEELS"
$ws.Range("F19").Value2 = "₹5.89 Lakh"
$ws.Range("G19").Value2 = "EMIs from ₹11,515/month"
$ws.Range("B20").Value2 = "2013 Maruti Wagon R 1.0 VXI"
$ws.Range("C20").Value2 = "48,692 KM"
$ws.Range("D20").Value2 = "MANUAL"
$ws.Range("E20").Value2 = "TOP MODEL"
$ws.Range("F20").Value2 = "₹2.79 Lakh"
$ws.Range("G20").Value2 = "EMIs from ₹7,347/month"
